$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.713.85'
$ws.Range("E2").Value = '  -1.44%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.025.49'
$ws.Range("E3").Value = '  -1.78%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.97'
$ws.Range("E5").Value = '  -1.71%  '

$ws.Range("E6").Value = '  -0.78%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.55'
$ws.Range("E7").Value = '  +1.32%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.383'
$ws.Range("E9").Value = '  -1.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0810'
$ws.Range("E10").Value = '  -0.03%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.103'
$ws.Range("E11").Value = '  -0.73%  '

$ws.Range("E12").Value = '  -1.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.324.24'
$ws.Range("E13").Value = '  -1.80%  '

$ws.Range("E14").Value = '  +0.97%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.755'
$ws.Range("E15").Value = '  -0.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.16'
$ws.Range("E16").Value = '  -2.63%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.031.22'
$ws.Range("E17").Value = '  -1.44%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.644.83'
$ws.Range("E18").Value = '  -1.28%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.02'
$ws.Range("E19").Value = '  -2.44%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.65'
$ws.Range("E20").Value = '  -0.45%  '

$ws.Range("E21").Value = '  -1.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '222.87'
$ws.Range("E22").Value = '  -1.04%  '

$ws.Range("E23").Value = '  +0.02%  '

$ws.Range("E25").Value = '  -2.57%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.23'
$ws.Range("E26").Value = '  -1.15%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '165.31'
$ws.Range("E27").Value = '  -0.60%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.129'
$ws.Range("E28").Value = '  -3.21%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.89'
$ws.Range("E29").Value = '  -1.11%  '

$ws.Range("E30").Value = '  -5.73%  '

$ws.Range("E31").Value = '  +0.76%  '

$ws.Range("E32").Value = '  -3.12%  '

$ws.Range("E33").Value = '  +4.23%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0600'
$ws.Range("E34").Value = '  -2.48%  '

$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.48'
$ws.Range("E35").Value = '  -3.20%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.32'
$ws.Range("E36").Value = '  +4.29%  '

$ws.Range("E37").Value = '  -3.95%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.24'
$ws.Range("E38").Value = '  -2.95%  '

$ws.Range("E39").Value = '  -0.20%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.535.39'
$ws.Range("E40").Value = '  +3.48%  '

$ws.Range("E41").Value = '  -1.45%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.68'
$ws.Range("E42").Value = '  -2.00%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.77'
$ws.Range("E43").Value = '  -1.02%  '

$ws.Range("E44").Value = '  -0.51%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0916'
$ws.Range("E45").Value = '  -3.67%  '

$ws.Range("E46").Value = '  -1.87%  '

$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.13'
$ws.Range("E47").Value = '  +0.45%  '

$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  -2.00%  '

$ws.Range("B49").Value = 'FTXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.87'
$ws.Range("E49").Value = '  -4.92%  '

$ws.Range("E50").Value = '  -0.56%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.215.53'
$ws.Range("E51").Value = '  -1.62%  '
